$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new column before E ("references") - shifts old E..J -> F..K
# ---------------------------------------------------------------------------
$ws.Columns("E:E").Insert()

# ---------------------------------------------------------------------------
# 2) Insert a new row before row 4 (new REQ006 item) - shifts old 4..6 -> 5..7
# ---------------------------------------------------------------------------
$ws.Rows("4:4").Insert()

# ---------------------------------------------------------------------------
# 3) Header row
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "uid"
$ws.Range("B1").Value = "level"
$ws.Range("C1").Value = "text"
$ws.Range("D1").Value = "ref"
$ws.Range("E1").Value = "references"
$ws.Range("F1").Value = "links"
$ws.Range("G1").Value = "active"
$ws.Range("H1").Value = "derived"
$ws.Range("I1").Value = "header"
$ws.Range("J1").Value = "normative"
$ws.Range("K1").Value = "reviewed"

# ---------------------------------------------------------------------------
# 4) Row 2 - REQ001 (existing row, data shifted right from col E)
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "SYS001`nSYS002:abc123"
$ws.Range("G2").Value = $true
$ws.Range("H2").Value = $false
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = $true
$ws.Range("K2").Value = ""

# ---------------------------------------------------------------------------
# 5) Row 3 - REQ003 (existing row, data shifted right from col E)
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "REQ001"
$ws.Range("G3").Value = $true
$ws.Range("H3").Value = $false
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = $true
$ws.Range("K3").Value = ""

# ---------------------------------------------------------------------------
# 6) Row 4 - NEW item REQ006 (brand new row)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "REQ006"
$ws.Range("B4").Value = "1.5"
$ws.Range("C4").Value = "Hello, world!"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "type:file,path:external/text.txt`ntype:file,path:external/text2.txt"
$ws.Range("F4").Value = "REQ001:35ed54323e3054c33ae5545fffdbbbf5"
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = $false
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = $true
$ws.Range("K4").Value = "c442316131ca0225595ae257f3b4583d"

# ---------------------------------------------------------------------------
# 7) Row 5 - REQ004 (previously row 4, now shifted to row 5)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "REQ004"
$ws.Range("B5").Value = "1.6"
$ws.Range("C5").Value = "Hello, world!"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = $true
$ws.Range("H5").Value = $false
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = $true
$ws.Range("K5").Value = ""

# ---------------------------------------------------------------------------
# 8) Row 6 - REQ002 (previously row 5, now shifted to row 6)
# ---------------------------------------------------------------------------
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = $false
$ws.Range("I6").Value = "Plantuml"
$ws.Range("J6").Value = $true
$ws.Range("K6").Value = "50ae164a198e612dee696cc80942dc29"

# ---------------------------------------------------------------------------
# 9) Row 7 - REQ2-001 (previously row 6, now shifted to row 7)
# ---------------------------------------------------------------------------
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = "REQ001"
$ws.Range("G7").Value = $true
$ws.Range("H7").Value = $false
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = $true
$ws.Range("K7").Value = ""

# ---------------------------------------------------------------------------
# 10) B4/B5 ("1.5"/"1.6") must stay text, not auto-convert to numbers.
#     Temporarily mark the cell as Text, re-enter the value, then copy the
#     (unformatted) style from a neighboring "s=2" cell back on top so the
#     cell ends up on the same style index as the rest of the row.
# ---------------------------------------------------------------------------
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "1.5"
$ws.Range("A4").Copy()
$ws.Range("B4").PasteSpecial(-4122) | Out-Null

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "1.6"
$ws.Range("A5").Copy()
$ws.Range("B5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 11) Column widths
#     Columns A-D and G-K already retain their correct widths after the
#     column insert above (the insert shifted old E..J -> F..K verbatim).
#     Only E (brand new) and F (changed from 16.5 -> 42.5) need fixing.
#     NOTE: Excel's ColumnWidth property is offset from the stored OOXML
#     <col width="..."> by exactly 5/6 (0.8333333333333333), so subtract
#     that constant to land on the exact target XML width.
# ---------------------------------------------------------------------------
$padding = 0.8333333333333333
$ws.Columns("E:E").ColumnWidth = 36.5 - $padding
$ws.Columns("F:F").ColumnWidth = 42.5 - $padding

# ---------------------------------------------------------------------------
# 12) Re-entering multi-line (wrapped) cell content above pins an explicit
#     row height on rows 2 and 4. AutoFit puts them back to the implicit
#     default height so no stray ht=/customHeight= survives in the XML.
# ---------------------------------------------------------------------------
$ws.Rows("1:7").AutoFit()

# ---------------------------------------------------------------------------
# 13) Fix up autofilter range (A1:J1 -> A1:K1)
# ---------------------------------------------------------------------------
if ($ws.AutoFilterMode) {
  $ws.AutoFilterMode = $false
}
$ws.Range("A1:K1").AutoFilter()

# ---------------------------------------------------------------------------
# 14) Fix up the hidden _FilterDatabase defined name
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
  if ($n.Name -like "*_FilterDatabase*") {
    $n.RefersTo = "=Sheet!`$A`$1:`$K`$1"
  }
}
